$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures scraped on 2023-01-15.
# Cells store values as literal text (matching the sheet's existing
# inline-string layout), so each target cell is formatted as Text
# before the new literal is written - this prevents Excel from
# auto-converting numeric-looking strings (e.g. "298.77", "-2.34%")
# into real numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.34%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.46%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.065"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07524"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.97%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.784"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.46%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.726"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "9.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.791"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.48%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9250"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.32%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1707"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07337"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.85%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07937"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.09%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03029"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.82%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.45%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.44%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04653"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.26%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006322"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.12%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.450"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.05%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.221"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.45%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.36%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.549"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.95%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1549"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.83%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.27%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004426"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.29%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.73%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001841"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-4.81%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01672"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04548"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.89%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007008"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.73%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.75%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002058"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.80%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01276"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006027"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7116"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-62.40%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-5.64%"
